$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.659021973609924
$ws.Range("B1").Value = 2.262895107269287
$ws.Range("C1").Value = 4.481653690338135
$ws.Range("D1").Value = 4.458285331726074
$ws.Range("E1").Value = 1.49140727519989
